$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting existing rows 20-50 down to 21-51
$ws.Rows("20:20").Insert()

# Populate the new row 20 with this week's data (values match the rest of the
# dataset's fixed columns; date/volume/price columns are the new observation)
$ws.Cells.Item(20, 1).Value = 8
$ws.Cells.Item(20, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 45133
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = 100112026
$ws.Cells.Item(20, 7).Value = "Haba"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 12000
$ws.Cells.Item(20, 12).Value = 13000
$ws.Cells.Item(20, 13).Value = 12500
$ws.Cells.Item(20, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 16).Value = 500
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
